# Applies the commit:
#   "custom accuracy + 데이터 1000개"
# which, for this sheet, means:
#   - the values in row 5 (B5:AH5) are rounded to a "custom accuracy" of
#     2 decimal places (from their original 3-decimal precision)
#   - the last data row (row 6) is removed entirely
#   - the sheet's dimension shrinks from A1:AH6 to A1:AH5 accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round-half-to-even on the exact decimal text of the double (avoids the
# binary floating point noise that a naive Math.Round would be sensitive
# to, and matches standard "banker's rounding" semantics).
function Round-HalfEven2 {
    param([double]$value)

    $neg = $false
    $v = $value
    if ($v -lt 0) {
        $neg = $true
        $v = -$v
    }

    $s = $v.ToString()
    if ($s.Contains("E")) {
        $s = $v.ToString("F10")
    }

    $dotIdx = $s.IndexOf(".")
    if ($dotIdx -lt 0) {
        if ($neg) { return (-1.0) * [double]$s } else { return [double]$s }
    }

    $intPart = $s.Substring(0, $dotIdx)
    $fracPart = $s.Substring($dotIdx + 1)
    $fracPart = $fracPart.PadRight(3, '0')

    $keep = $fracPart.Substring(0, 2)
    $rest = $fracPart.Substring(2)

    $combinedStr = [string]::Concat($intPart, $keep)
    $combined = [long]$combinedStr

    $isExactHalf = $false
    if ($rest.Length -gt 0 -and $rest.Substring(0, 1) -eq '5') {
        $restTail = $rest.Substring(1)
        $restTail = $restTail.TrimEnd('0')
        if ($restTail.Length -eq 0) {
            $isExactHalf = $true
        }
    }

    $restVal = [double]([string]::Concat("0.", $rest))

    $roundUp = $false
    if ($restVal -gt 0.5) {
        $roundUp = $true
    } elseif ($isExactHalf) {
        $lastDigitChar = $keep.Substring(1, 1)
        $lastDigit = [int]$lastDigitChar
        if (($lastDigit % 2) -ne 0) {
            $roundUp = $true
        }
    }

    if ($roundUp) {
        $combined = $combined + 1
    }

    $combinedStr2 = $combined.ToString()
    $combinedStr2 = $combinedStr2.PadLeft(3, '0')
    $len = $combinedStr2.Length
    $newIntPart = $combinedStr2.Substring(0, $len - 2)
    $newFracPart = $combinedStr2.Substring($len - 2)
    $resultStr = [string]::Concat($newIntPart, ".", $newFracPart)

    $result = [double]$resultStr
    if ($neg) { $result = -1.0 * $result }
    return $result
}

# Round B5:AH5 (columns 2..34) down to 2 decimals of accuracy.
for ($c = 2; $c -le 34; $c++) {
    $cell = $ws.Cells.Item(5, $c)
    $orig = [double]$cell.Value2
    $cell.Value = Round-HalfEven2 $orig
}

# Drop the last data row (row 6) entirely; the sheet dimension is
# recalculated automatically (A1:AH6 -> A1:AH5).
$ws.Rows.Item(6).Delete()
